$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Shape 1 (title): add "Exemplo Loja" ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$titleRun = $title.InsertAfter("Exemplo Loja")
$titleRun.LanguageID = "pt-BR"

# --- Shape 2 (content placeholder): add three paragraphs of text ---
$body = $s.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1
$r = $body.InsertAfter(" Classes Monitor ")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32

$r = $body.InsertAfter("e Impressora.")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32

$r = $body.InsertAfter("`r")

# Paragraph 2
$r = $body.InsertAfter(" Novos métodos são implementados nas classes correspondentes ")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32

$r = $body.InsertAfter("VisitorPeso")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32

$r = $body.InsertAfter(" e ")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32

$r = $body.InsertAfter("VisitorPreco")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32

$r = $body.InsertAfter(".")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32

$r = $body.InsertAfter("`r")

# Paragraph 3
$r = $body.InsertAfter(" ")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32

$r = $body.InsertAfter("As classes originais permanecem inalteradas.")
$r.LanguageID = "pt-BR"
$r.Font.Size = 32
